# The sheet originally had a header row (model_number / parts / total_bays)
# in row 1, with the actual data rows (MN1/MN2) below it. The export to
# Postgres needs the raw data only, so drop the header row entirely and
# let the data shift up.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("1").Delete()

# Leave the sheet's selection on the (now-header-less) first row, matching
# a "select row 1" state rather than pointing at a stray cell below the data.
$ws.Rows("1").Select() | Out-Null
